$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; this shifts the existing rows 40-84
# down to 41-85 (and inherits the row-40 formatting, e.g. the date style
# on column D).
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new price record.
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = 'Vega Modelo de Temuco'
$ws.Range("C40").Value = 'La Araucanía'
$ws.Range("D40").Value = 44587
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 'Fruta'
$ws.Range("G40").Value = 100101
$ws.Range("H40").Value = 'Berries'
$ws.Range("I40").Value = 100101001
$ws.Range("J40").Value = 'Arándano (blue)'
$ws.Range("K40").Value = 'Sin especificar'
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 100
$ws.Range("N40").Value = 2000
$ws.Range("O40").Value = 2000
$ws.Range("P40").Value = 2000
$ws.Range("Q40").Value = '$/kilo'
$ws.Range("R40").Value = 'Región del Maule'
$ws.Range("S40").Value = 2000
$ws.Range("T40").Value = 1
